$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 147; this shifts existing rows 147..261 down to 148..262.
$ws.Rows.Item(147).Insert()

# Fill in the new row 147 with the new record's data.
$ws.Cells.Item(147, 1).Value = 5
$ws.Cells.Item(147, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(147, 3).Value = "Maule"
$ws.Cells.Item(147, 4).Value = 44574
$ws.Cells.Item(147, 5).Value = 7
$ws.Cells.Item(147, 6).Value = 100114013
$ws.Cells.Item(147, 7).Value = "Zanahoria"
$ws.Cells.Item(147, 8).Value = "Sin especificar"
$ws.Cells.Item(147, 9).Value = "Primera"
$ws.Cells.Item(147, 10).Value = 400
$ws.Cells.Item(147, 11).Value = 7000
$ws.Cells.Item(147, 12).Value = 7000
$ws.Cells.Item(147, 13).Value = 7000
$ws.Cells.Item(147, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(147, 15).Value = "Región de Ñuble"
$ws.Cells.Item(147, 16).Value = 350
$ws.Cells.Item(147, 17).Value = 20
$ws.Cells.Item(147, 18).Value = "Hortaliza"

# Match the date style used by the rest of column D.
$ws.Cells.Item(147, 4).NumberFormat = $ws.Cells.Item(148, 4).NumberFormat
